# Apply the StructureDefinition-region-code.xlsx update:
#  - Metadata sheet: URL / Version / Date / Publisher refreshed for the
#    LinuxForHealth re-brand (was Alvearie / ibm.com).
#  - Elements sheet: the root "Extension" row's rolled-up Constraint(s)
#    text (ele-1/ext-1) is cleared - it now only lives on the
#    Extension.extension child row.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/region-code"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
